$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H4").Value = 3790.75
$ws.Range("I4").Value = 2248.9
$ws.Range("K4").Value = 2248.9
$ws.Range("M4").Value = -2134.9
$ws.Range("H15").Value = 91.29412000000001
$ws.Range("I15").Value = 91.29412000000001
$ws.Range("K15").Value = 273.88236
$ws.Range("M15").Value = -104.88236
$ws.Range("H38").Value = 338
$ws.Range("I38").Value = 155.25
$ws.Range("J38").Value = 1800
$ws.Range("K38").Value = 465.75
$ws.Range("L38").Value = 5400
$ws.Range("M38").Value = -93.75
$ws.Range("N38").Value = -6144
$ws.Range("H43").Value = 3299.1428
$ws.Range("J43").Value = 3266.6667
$ws.Range("L43").Value = 3266.6667
$ws.Range("N43").Value = -3404.6667
$ws.Range("H46").Value = 10000
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 30000
$ws.Range("N46").Value = -30238
$ws.Range("H60").Value = 10000
$ws.Range("J60").Value = 10000
$ws.Range("L60").Value = 30000
$ws.Range("N60").Value = -30968
$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 3000
$ws.Range("N64").Value = 0
$ws.Range("M64").Value = -2752
$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 3000
$ws.Range("N67").Value = 0
$ws.Range("M67").Value = -2142
$ws.Range("H92").Value = 528.4
$ws.Range("I92").Value = 537.1111
$ws.Range("J92").Value = 450
$ws.Range("K92").Value = 537.1111
$ws.Range("L92").Value = 450
$ws.Range("M92").Value = 710.8889
$ws.Range("N92").Value = -2946
$ws.Range("H106").Value = 2488.5
$ws.Range("I106").Value = 1721.25
$ws.Range("K106").Value = 1721.25
$ws.Range("M106").Value = -1090.25
$ws.Range("H110").Value = 79997
$ws.Range("I110").Value = 70000
$ws.Range("J110").Value = 89994
$ws.Range("K110").Value = 70000
$ws.Range("L110").Value = 89994
$ws.Range("M110").Value = -65910
$ws.Range("N110").Value = -98174
$ws.Range("H125").Value = 730.4
$ws.Range("I125").Value = 713
$ws.Range("K125").Value = 6417
$ws.Range("M125").Value = -3957
$ws.Range("H132").Value = 2649.0386
$ws.Range("I132").Value = 2649.0386
$ws.Range("K132").Value = 7947.1158
$ws.Range("M132").Value = -5417.1158
$ws.Range("H137").Value = 2000
$ws.Range("J137").Value = 2000
$ws.Range("L137").Value = 6000
$ws.Range("N137").Value = -11100
$ws.Range("L64").ClearContents()
$ws.Range("L67").ClearContents()

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("H74").Value = 1649.75
$ws.Range("I74").Value = 1649.75
$ws.Range("K74").Value = 1649.75
$ws.Range("M74").Value = -775.75
$ws.Range("H77").Value = 1649.75
$ws.Range("I77").Value = 1649.75
$ws.Range("K77").Value = 8248.75
$ws.Range("M77").Value = -3880.75
$ws.Range("H110").Value = 941.9091
$ws.Range("I110").Value = 941.9091
$ws.Range("K110").Value = 941.9091
$ws.Range("M110").Value = 1103.0909
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("M136").ClearContents()

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H86").Value = 25750.75
$ws.Range("J86").Value = 25667.834
$ws.Range("L86").Value = 25667.834
$ws.Range("N86").Value = -27913.834
$ws.Range("H89").Value = 25750.75
$ws.Range("J89").Value = 25667.834
$ws.Range("L89").Value = 128339.17
$ws.Range("N89").Value = -139571.17
$ws.Range("H94").Value = 3541.3076
$ws.Range("I94").Value = 2375.182
$ws.Range("K94").Value = 2375.182
$ws.Range("M94").Value = -1924.182

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H15").Value = 5102.6665
$ws.Range("J15").Value = 7504
$ws.Range("L15").Value = 7504
$ws.Range("N15").Value = -7844
$ws.Range("H16").Value = 799
$ws.Range("I16").Value = 799
$ws.Range("K16").Value = 799
$ws.Range("M16").Value = -512
$ws.Range("H105").Value = 3166.3333
$ws.Range("I105").Value = 999
$ws.Range("J105").Value = 4250
$ws.Range("K105").Value = 999
$ws.Range("L105").Value = 4250
$ws.Range("M105").Value = 748
$ws.Range("N105").Value = -7744
$ws.Range("H107").Value = 649.1875
$ws.Range("I107").Value = 484.85715
$ws.Range("K107").Value = 484.85715
$ws.Range("M107").Value = 1435.14285
$ws.Range("H113").Value = 799
$ws.Range("I113").Value = 799
$ws.Range("K113").Value = 799
$ws.Range("M113").Value = 1371

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H4").Value = 316145
$ws.Range("I4").Value = 294485.66
$ws.Range("K4").Value = 883456.98
$ws.Range("M4").Value = -883344.98
$ws.Range("H39").Value = 3954.4666
$ws.Range("I39").Value = 2536.6667
$ws.Range("J39").Value = 4308.9165
$ws.Range("K39").Value = 7610.000100000001
$ws.Range("L39").Value = 12926.7495
$ws.Range("M39").Value = -7316.000100000001
$ws.Range("N39").Value = -13514.7495
$ws.Range("H55").Value = 2984.5833
$ws.Range("I55").Value = 1900
$ws.Range("J55").Value = 3201.5
$ws.Range("K55").Value = 5700
$ws.Range("L55").Value = 9604.5
$ws.Range("M55").Value = -5523
$ws.Range("N55").Value = -9958.5
$ws.Range("H62").Value = 19000
$ws.Range("I62").Value = 19000
$ws.Range("K62").Value = 57000
$ws.Range("M62").Value = -56314
$ws.Range("H65").Value = 19000
$ws.Range("I65").Value = 19000
$ws.Range("K65").Value = 171000
$ws.Range("M65").Value = -167568

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H5").Value = 18999.5
$ws.Range("I5").Value = 18999
$ws.Range("J5").Value = 19000
$ws.Range("K5").Value = 18999
$ws.Range("L5").Value = 19000
$ws.Range("M5").Value = -18887
$ws.Range("N5").Value = -19224
$ws.Range("H113").Value = 1258
$ws.Range("I113").Value = 1197.5
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 1197.5
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = 972.5
$ws.Range("N113").Value = -5840
$ws.Range("H114").Value = 99995
$ws.Range("J114").Value = 99995
$ws.Range("L114").Value = 99995
$ws.Range("N114").Value = -108673

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H7").Value = 6098.2
$ws.Range("I7").Value = 5622.75
$ws.Range("J7").Value = 8000
$ws.Range("K7").Value = 5622.75
$ws.Range("L7").Value = 8000
$ws.Range("M7").Value = -5510.75
$ws.Range("N7").Value = -8224
$ws.Range("H80").Value = 55499.5
$ws.Range("I80").Value = 45999
$ws.Range("J80").Value = 65000
$ws.Range("K80").Value = 45999
$ws.Range("L80").Value = 65000
$ws.Range("M80").Value = -44876
$ws.Range("N80").Value = -67246
$ws.Range("H83").Value = 55499.5
$ws.Range("I83").Value = 45999
$ws.Range("J83").Value = 65000
$ws.Range("K83").Value = 137997
$ws.Range("L83").Value = 195000
$ws.Range("M83").Value = -132381
$ws.Range("N83").Value = -206232
$ws.Range("H100").Value = 11250
$ws.Range("I100").Value = 11250
$ws.Range("K100").Value = 11250
$ws.Range("M100").Value = -10709
$ws.Range("H122").Value = 5166.3335
$ws.Range("I122").Value = 5166.3335
$ws.Range("K122").Value = 15499.0005
$ws.Range("M122").Value = -13049.0005
$ws.Range("H126").Value = 6098.2
$ws.Range("I126").Value = 5622.75
$ws.Range("J126").Value = 8000
$ws.Range("K126").Value = 16868.25
$ws.Range("L126").Value = 24000
$ws.Range("M126").Value = -14398.25
$ws.Range("N126").Value = -28940

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H2").Value = 5125000
$ws.Range("I2").Value = 5125000
$ws.Range("K2").Value = 5125000
$ws.Range("M2").Value = -5124888
$ws.Range("H4").Value = 494232.1
$ws.Range("I4").Value = 80930.39999999999
$ws.Range("J4").Value = 838650.2
$ws.Range("K4").Value = 80930.39999999999
$ws.Range("L4").Value = 838650.2
$ws.Range("M4").Value = -80817.39999999999
$ws.Range("N4").Value = -838876.2
$ws.Range("H100").Value = 15449.8
$ws.Range("I100").Value = 37649.5
$ws.Range("K100").Value = 75299
$ws.Range("M100").Value = -74758
$ws.Range("H126").Value = 3709.6365
$ws.Range("I126").Value = 2600.875
$ws.Range("J126").Value = 6666.3335
$ws.Range("K126").Value = 7802.625
$ws.Range("L126").Value = 19999.0005
$ws.Range("M126").Value = -5332.625
$ws.Range("N126").Value = -24939.0005
